{"js": "// Commit: \"fixed discounting factor mistake\"\n// The sentence incorrectly said delta was the probability of the repeated\n// game ENDING; it should say delta is the probability of the game NOT\n// ending (discounting interpretation).\nconst oldText = \"as the probability of the repeated game ending then the\";\nconst newText = \"as the probability of the repeated game not ending then the\";\n\nconst searchResults = context.document.body.search(oldText, {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the sentence to fix.\");\n}\n\nsearchResults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Commit: \"fixed discounting factor mistake\"\n# The sentence incorrectly said delta was the probability of the repeated\n# game ENDING; it should say delta is the probability of the game NOT\n# ending (discounting interpretation).\n\n$d = $word.ActiveDocument\n\n$oldText = \"as the probability of the repeated game ending then the\"\n$newText = \"as the probability of the repeated game not ending then the\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceOne (1): replace only the first occurrence (the string is unique).\n$found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 1)\nif (-not $found) {\n    throw \"Could not find the sentence to fix.\"\n}\n"}
